$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Checklist" -> "Session"
$ws.Name = "Session"

# Update "Selection" -> "Scan" for rows 2-5 (column E)
$ws.Range("E2").Value = "Scan"
$ws.Range("E3").Value = "Scan"
$ws.Range("E4").Value = "Scan"
$ws.Range("E5").Value = "Scan"

# Row 6 gets the values that used to be on row 7's A/D columns
$ws.Range("A6").Value = "'555585"
$ws.Range("A6").ClearFormats()
$ws.Range("D6").Value = "12:44:47"

# Delete the old row 7 (shifts rows up, removing trailing row)
$ws.Rows.Item(7).Delete()
